$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting the existing row 79 (and everything
# below it) down by one. This grows the used range from A1:R189 to A1:R190.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with its data (a new weekly price entry).
$ws.Cells.Item(79, 1).Value = 8
$ws.Cells.Item(79, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(79, 3).Value = "Coquimbo"
$ws.Cells.Item(79, 4).Value = 44915
$ws.Cells.Item(79, 5).Value = 4
$ws.Cells.Item(79, 6).Value = 100112040
$ws.Cells.Item(79, 7).Value = "Cilantro"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 2400
$ws.Cells.Item(79, 11).Value = 2500
$ws.Cells.Item(79, 12).Value = 3000
$ws.Cells.Item(79, 13).Value = 2750
$ws.Cells.Item(79, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(79, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(79, 16).Value = 1833
$ws.Cells.Item(79, 17).Value = 1.5
$ws.Cells.Item(79, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(79, 4).NumberFormat = $ws.Cells.Item(80, 4).NumberFormat
